# Assignment 2 test-plan update: fill in the SavingsAccount unit-test rows,
# update the Developer cell, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name (was "Student Name" placeholder)
$ws.Range("C3").Value = "Md Apurba Khan"

# --- Test case 1 (row 7): __init__ / valid attributes ---
$ws.Range("E7").Value = "Valid account_number, client_number, balance, date_created, minimum_balance"
$ws.Range("F7").Value = "SavingsAccount(5001, 1003, 150, date(2022, 5, 10), 50.0)"
$ws.Range("G7").Value = "Instance is created successfully with correct attributes."

# --- Test case 2 (row 8): __init__ / invalid minimum_balance type ---
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'SavingsAccount(5002, 1004, 200, date(2023, 6, 10), "invalid")'
$ws.Range("G8").Value = "minimum_balance defaults to 50.0."

# --- Test case 3 (row 9): get_service_charges / balance greater than minimum ---
$ws.Range("E9").Value = "self._balance = 200, self._minimum_balance = 50"
$ws.Range("F9").Value = "get_service_charges()"
$ws.Range("G9").Value = "Returns base service charge (e.g., `$0.50)."

# --- Test case 4 (row 10): get_service_charges / balance equal to minimum ---
$ws.Range("E10").Value = "self._balance = 50, self._minimum_balance = 50"
$ws.Range("F10").Value = "get_service_charges()"
$ws.Range("G10").Value = "Returns base service charge (e.g., `$0.50)."

# --- Test case 5 (row 11): get_service_charges / balance less than minimum ---
$ws.Range("E11").Value = "self._balance = 20, self._minimum_balance = 50"
$ws.Range("F11").Value = "get_service_charges()"
$ws.Range("G11").Value = "Returns BASE_SERVICE_CHARGE * SERVICE_CHARGE_PREMIUM (e.g., `$1.00)."

# --- Test case 6 (row 12): __str__ / correct formatted string ---
$ws.Range("E12").Value = "Instance has valid attributes set"
$ws.Range("F12").Value = "str(savings_account)"
$ws.Range("G12").Value = "Returns a formatted string including account number, balance, and minimum balance."

# Move the active selection to G12 (single cell), replacing the old B7:B32 selection
[void]$ws.Range("G12").Select()
